$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns so numeric-looking strings
# (e.g. "244.72", "0.630") are preserved as text instead of being auto-converted
# to numbers by Excel's smart-entry parsing.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '41.292.42'
$ws.Range("E2").Value = '  -5.88%  '

$ws.Range("D3").Value = '2.223.24'
$ws.Range("E3").Value = '  -5.98%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '244.72'
$ws.Range("E5").Value = '  +2.04%  '

$ws.Range("D6").Value = '0.630'
$ws.Range("E6").Value = '  -6.04%  '

$ws.Range("D7").Value = '69.48'
$ws.Range("E7").Value = '  -6.33%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").Value = '0.559'
$ws.Range("E9").Value = '  -6.61%  '

$ws.Range("D10").Value = '39.34'
$ws.Range("E10").Value = '  +5.86%  '

$ws.Range("D11").Value = '0.0957'
$ws.Range("E11").Value = '  -6.82%  '

$ws.Range("D12").Value = '58.27'
$ws.Range("E12").Value = '  -3.14%  '

$ws.Range("E13").Value = '  -4.31%  '

$ws.Range("D14").Value = '6.79'
$ws.Range("E14").Value = '  -6.74%  '

$ws.Range("D15").Value = '2.553.29'
$ws.Range("E15").Value = '  -6.09%  '

$ws.Range("D16").Value = '14.85'
$ws.Range("E16").Value = '  -9.65%  '

$ws.Range("D17").Value = '0.841'
$ws.Range("E17").Value = '  -9.33%  '

$ws.Range("D18").Value = '2.218.32'
$ws.Range("E18").Value = '  -6.45%  '

$ws.Range("D19").Value = '41.317.15'
$ws.Range("E19").Value = '  -5.72%  '

$ws.Range("D20").Value = '0.0₃0956'
$ws.Range("E20").Value = '  -8.04%  '

$ws.Range("D21").Value = '72.34'
$ws.Range("E21").Value = '  -6.41%  '

$ws.Range("D22").Value = '6.10'
$ws.Range("E22").Value = '  -7.70%  '

$ws.Range("D23").Value = '232.39'
$ws.Range("E23").Value = '  -8.49%  '

$ws.Range("E24").Value = '  +13.31%  '

$ws.Range("E25").Value = '  +0.14%  '

$ws.Range("E26").Value = '  -4.57%  '

$ws.Range("E27").Value = '  -2.94%  '

$ws.Range("D28").Value = '9.82'
$ws.Range("E28").Value = '  -7.65%  '

$ws.Range("E29").Value = '  -4.91%  '

$ws.Range("D30").Value = '171.49'
$ws.Range("E30").Value = '  -2.07%  '

$ws.Range("D31").Value = '20.58'
$ws.Range("E31").Value = '  -8.08%  '

$ws.Range("E32").Value = '  -7.55%  '

$ws.Range("D34").Value = '0.0717'
$ws.Range("E34").Value = '  -5.10%  '

$ws.Range("E35").Value = '  -4.03%  '

$ws.Range("D36").Value = '4.62'
$ws.Range("E36").Value = '  -9.71%  '

$ws.Range("D37").Value = '3.91'
$ws.Range("E37").Value = '  +2.66%  '

$ws.Range("E38").Value = '  +16.32%  '

$ws.Range("E39").Value = '  -0.77%  '

$ws.Range("E40").Value = '  -4.73%  '

$ws.Range("D41").Value = '5.86'
$ws.Range("E41").Value = '  -11.49%  '

$ws.Range("D42").Value = '66.36'
$ws.Range("E42").Value = '  +2.13%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '0.208'
$ws.Range("E43").Value = '  +3.35%  '

$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").Value = '5.00'
$ws.Range("E44").Value = '  -11.56%  '

$ws.Range("D45").Value = '8.90'
$ws.Range("E45").Value = '  -1.78%  '

$ws.Range("E46").Value = '  -6.79%  '

$ws.Range("D47").Value = '10.82'
$ws.Range("E47").Value = '  +11.34%  '

$ws.Range("D48").Value = '4.64'
$ws.Range("E48").Value = '  +6.25%  '

$ws.Range("E49").Value = '  +0.05%  '

$ws.Range("E50").Value = '  -5.70%  '

$ws.Range("E51").Value = '  -5.31%  '

# Restore the default cell style so no stray formatting is introduced
# (matches the original workbook, where these cells carry no explicit style).
$priceRange.Style = "Normal"
